$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 810.7586
$ws.Range("I28").Value = 428.43478
$ws.Range("K28").Value = 428.43478
$ws.Range("M28").Value = 56.56522000000001

$ws.Range("H32").Value = 11378.25
$ws.Range("J32").Value = 11547.75
$ws.Range("L32").Value = 11547.75
$ws.Range("N32").Value = -12199.75

$ws.Range("H74").Value = 8500
$ws.Range("I74").Value = 8500
$ws.Range("K74").Value = 8500
$ws.Range("M74").Value = -7564

$ws.Range("H77").Value = 8500
$ws.Range("I77").Value = 8500
$ws.Range("K77").Value = 42500
$ws.Range("M77").Value = -37820

$ws.Range("H101").Value = 314.55554
$ws.Range("J101").Value = 200
$ws.Range("L101").Value = 600
$ws.Range("N101").Value = -3844

$ws.Range("H103").Value = 747.8
$ws.Range("I103").Value = 520.6667
$ws.Range("J103").Value = 899.2222
$ws.Range("K103").Value = 1562.0001
$ws.Range("L103").Value = 2697.6666
$ws.Range("M103").Value = -976.0001
$ws.Range("N103").Value = -3869.6666

$ws.Range("H135").Value = 1107.65
$ws.Range("I135").Value = 509
$ws.Range("J135").Value = 4500
$ws.Range("K135").Value = 4581
$ws.Range("L135").Value = 40500
$ws.Range("M135").Value = -2046
$ws.Range("N135").Value = -45570

$ws.Range("H137").Value = 19338.223
$ws.Range("I137").Value = 23999.572
$ws.Range("K137").Value = 71998.716
$ws.Range("M137").Value = -69448.716

$ws.Range("H138").Value = 35015.16
$ws.Range("I138").Value = 2077.4443
$ws.Range("K138").Value = 6232.3329
$ws.Range("M138").Value = -1092.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29310.076
$ws.Range("I32").Value = 31474.805
$ws.Range("J32").Value = 3333.3333
$ws.Range("K32").Value = 31474.805
$ws.Range("L32").Value = 3333.3333
$ws.Range("M32").Value = -31187.805
$ws.Range("N32").Value = -3907.3333

$ws.Range("H61").Value = 5616.893
$ws.Range("I61").Value = 1159.6666
$ws.Range("J61").Value = 13639.9
$ws.Range("K61").Value = 1159.6666
$ws.Range("L61").Value = 13639.9
$ws.Range("M61").Value = -947.6666
$ws.Range("N61").Value = -14063.9

$ws.Range("H63").Value = 2990
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 3112.5
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 3112.5
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -4484.5

$ws.Range("H66").Value = 2990
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 3112.5
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 15562.5
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -22426.5

$ws.Range("H132").Value = 1555.8793
$ws.Range("I132").Value = 994.19446
$ws.Range("K132").Value = 2982.58338
$ws.Range("M132").Value = -452.58338

$ws.Range("H136").Value = 5616.893
$ws.Range("I136").Value = 1159.6666
$ws.Range("J136").Value = 13639.9
$ws.Range("K136").Value = 3478.9998
$ws.Range("L136").Value = 40919.7
$ws.Range("M136").Value = -928.9998000000001
$ws.Range("N136").Value = -46019.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1714
$ws.Range("I86").Value = 1622.5834
$ws.Range("K86").Value = 1622.5834
$ws.Range("M86").Value = -499.5834

$ws.Range("H89").Value = 1714
$ws.Range("I89").Value = 1622.5834
$ws.Range("K89").Value = 8112.916999999999
$ws.Range("M89").Value = -2496.916999999999

$ws.Range("H105").Value = 4264.222
$ws.Range("I105").Value = 4066.7144
$ws.Range("K105").Value = 4066.7144
$ws.Range("M105").Value = -2319.7144

$ws.Range("H134").Value = 1910.8096
$ws.Range("I134").Value = 1257.2812
$ws.Range("J134").Value = 4002.1
$ws.Range("K134").Value = 3771.8436
$ws.Range("L134").Value = 12006.3
$ws.Range("M134").Value = -1236.8436
$ws.Range("N134").Value = -17076.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 834.8
$ws.Range("I22").Value = 699.75
$ws.Range("J22").Value = 924.8333
$ws.Range("K22").Value = 699.75
$ws.Range("L22").Value = 924.8333
$ws.Range("M22").Value = -349.75
$ws.Range("N22").Value = -1624.8333

$ws.Range("H31").Value = 5884789
$ws.Range("I31").Value = 7693415.5
$ws.Range("J31").Value = 6753
$ws.Range("K31").Value = 7693415.5
$ws.Range("L31").Value = 6753
$ws.Range("M31").Value = -7693120.5
$ws.Range("N31").Value = -7343

$ws.Range("H34").Value = 5884789
$ws.Range("I34").Value = 7693415.5
$ws.Range("J34").Value = 6753
$ws.Range("K34").Value = 7693415.5
$ws.Range("L34").Value = 6753
$ws.Range("M34").Value = -7693213.5
$ws.Range("N34").Value = -7157

$ws.Range("H58").Value = 17349.572
$ws.Range("I58").Value = 1770.7646
$ws.Range("K58").Value = 1770.7646
$ws.Range("M58").Value = -1567.7646

$ws.Range("H94").Value = 1135.4783
$ws.Range("I94").Value = 957.44446
$ws.Range("J94").Value = 1249.9286
$ws.Range("K94").Value = 957.44446
$ws.Range("L94").Value = 1249.9286
$ws.Range("M94").Value = -506.44446
$ws.Range("N94").Value = -2151.9286

$ws.Range("H105").Value = 25123.688
$ws.Range("I105").Value = 32900.332
$ws.Range("J105").Value = 1793.75
$ws.Range("K105").Value = 32900.332
$ws.Range("L105").Value = 1793.75
$ws.Range("M105").Value = -31153.332
$ws.Range("N105").Value = -5287.75

$ws.Range("H132").Value = 33308.258
$ws.Range("I132").Value = 39352.23
$ws.Range("K132").Value = 118056.69
$ws.Range("M132").Value = -115526.69

$ws.Range("H134").Value = 1494.5217
$ws.Range("I134").Value = 1335.1818
$ws.Range("K134").Value = 4005.5454
$ws.Range("M134").Value = -1470.5454

$ws.Range("H136").Value = 17349.572
$ws.Range("I136").Value = 1770.7646
$ws.Range("K136").Value = 5312.293799999999
$ws.Range("M136").Value = -2762.293799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 243.88889
$ws.Range("J12").Value = 163.57143
$ws.Range("L12").Value = 490.71429
$ws.Range("N12").Value = -836.71429

$ws.Range("H33").Value = 121.5
$ws.Range("I33").Value = 105.8
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 634.8
$ws.Range("L33").Value = 1200
$ws.Range("M33").Value = -351.8
$ws.Range("N33").Value = -1766

$ws.Range("H40").Value = 31
$ws.Range("I40").Value = 31
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 124
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -55
$ws.Range("N40").ClearContents()

$ws.Range("H44").Value = 10667
$ws.Range("I44").Value = 40003
$ws.Range("J44").Value = 4799.8
$ws.Range("K44").Value = 120009
$ws.Range("L44").Value = 14399.4
$ws.Range("M44").Value = -119611
$ws.Range("N44").Value = -15195.4

$ws.Range("H69").Value = 3757.1428
$ws.Range("J69").Value = 4281.8184
$ws.Range("L69").Value = 12845.4552
$ws.Range("N69").Value = -14467.4552

$ws.Range("H72").Value = 3757.1428
$ws.Range("J72").Value = 4281.8184
$ws.Range("L72").Value = 38536.3656
$ws.Range("N72").Value = -46648.3656

$ws.Range("H113").Value = 765.4
$ws.Range("J113").Value = 874.75
$ws.Range("L113").Value = 2624.25
$ws.Range("N113").Value = -6964.25

$ws.Range("H122").Value = 918.88464
$ws.Range("I122").Value = 285.625
$ws.Range("J122").Value = 1200.3334
$ws.Range("K122").Value = 2570.625
$ws.Range("L122").Value = 10803.0006
$ws.Range("M122").Value = -120.625
$ws.Range("N122").Value = -15703.0006

$ws.Range("H131").Value = 1590.5758
$ws.Range("J131").Value = 1975.381
$ws.Range("L131").Value = 5926.143
$ws.Range("N131").Value = -16006.143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6824.154
$ws.Range("I80").Value = 1919.75
$ws.Range("K80").Value = 1919.75
$ws.Range("M80").Value = -921.75

$ws.Range("H83").Value = 6824.154
$ws.Range("I83").Value = 1919.75
$ws.Range("K83").Value = 9598.75
$ws.Range("M83").Value = -4606.75

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H102").Value = 45807.145
$ws.Range("I102").Value = 69156.336
$ws.Range("J102").Value = 3778.6
$ws.Range("K102").Value = 69156.336
$ws.Range("L102").Value = 3778.6
$ws.Range("M102").Value = -67534.336
$ws.Range("N102").Value = -7022.6

$ws.Range("H126").Value = 2672.842
$ws.Range("I126").Value = 2056.4614
$ws.Range("J126").Value = 4008.3333
$ws.Range("K126").Value = 6169.3842
$ws.Range("L126").Value = 12024.9999
$ws.Range("M126").Value = -3699.3842
$ws.Range("N126").Value = -16964.9999

$ws.Range("H132").Value = 2877.1667
$ws.Range("I132").Value = 2849.375
$ws.Range("J132").Value = 3099.5
$ws.Range("K132").Value = 8548.125
$ws.Range("L132").Value = 9298.5
$ws.Range("M132").Value = -6018.125
$ws.Range("N132").Value = -14358.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1589.7646
$ws.Range("I22").Value = 1575.6666
$ws.Range("J22").Value = 1605.625
$ws.Range("K22").Value = 1575.6666
$ws.Range("L22").Value = 1605.625
$ws.Range("M22").Value = -1280.6666
$ws.Range("N22").Value = -2195.625

$ws.Range("H27").Value = 1589.7646
$ws.Range("I27").Value = 1575.6666
$ws.Range("J27").Value = 1605.625
$ws.Range("K27").Value = 1575.6666
$ws.Range("L27").Value = 1605.625
$ws.Range("M27").Value = -1468.6666
$ws.Range("N27").Value = -1819.625

$ws.Range("H132").Value = 1805.6666
$ws.Range("I132").Value = 899.9
$ws.Range("J132").Value = 2937.875
$ws.Range("K132").Value = 2699.7
$ws.Range("L132").Value = 8813.625
$ws.Range("M132").Value = -169.6999999999998
$ws.Range("N132").Value = -13873.625

$ws.Range("H136").Value = 3769.8333
$ws.Range("I136").Value = 3803.5293
$ws.Range("J136").Value = 3688
$ws.Range("K136").Value = 11410.5879
$ws.Range("L136").Value = 11064
$ws.Range("M136").Value = -8860.5879
$ws.Range("N136").Value = -16164
